$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update stock quantities (column B) for the restock datatable
$ws.Range("B2").Value = 8
$ws.Range("B3").Value = 29
$ws.Range("B5").Value = 1
$ws.Range("B6").Value = 129

# Update the active selection to B4
$ws.Range("B4").Select()
